# Regenerate the "K" column (column G) values on Sheet1.
# Commit intent: "regen save_data to use K instead of Strike#, regen
# std/mean, calc and write s_vals" -- the simulation/back-end recomputed
# the K statistic for every game row and rewrote column G in place.
# Column G is numeric ("K") in both before/after states; only the values
# change, row-by-row, for rows 2..64 (with rows 54 and 62 unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K value, taken from the regenerated save data.
$newK = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 2
    6  = 1
    7  = 0
    8  = 4
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 1
    14 = 2
    15 = 2
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 1
    22 = 1
    23 = 0
    24 = 2
    25 = 1
    26 = 0
    27 = 0
    28 = 0
    29 = 1
    30 = 0
    31 = 1
    32 = 1
    33 = 1
    34 = 0
    35 = 0
    36 = 0
    37 = 0
    38 = 0
    39 = 0
    40 = 2
    41 = 3
    42 = 0
    43 = 0
    44 = 0
    45 = 3
    46 = 1
    47 = 2
    48 = 0
    49 = 1
    50 = 2
    51 = 2
    52 = 2
    53 = 2
    55 = 1
    56 = 1
    57 = 1
    58 = 2
    59 = 2
    60 = 2
    61 = 4
    63 = 1
    64 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
